$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 787
$ws.Cells.Item(3, 6).Value = 519
$ws.Cells.Item(4, 6).Value = 265
$ws.Cells.Item(5, 6).Value = 478
$ws.Cells.Item(6, 6).Value = 1117
$ws.Cells.Item(7, 6).Value = 316
$ws.Cells.Item(8, 6).Value = 22
$ws.Cells.Item(11, 6).Value = 1108
$ws.Cells.Item(14, 6).Value = 755
$ws.Cells.Item(15, 6).Value = 795
$ws.Cells.Item(16, 6).Value = 177
$ws.Cells.Item(18, 6).Value = 59
$ws.Cells.Item(19, 6).Value = 661
$ws.Cells.Item(20, 6).Value = 158
$ws.Cells.Item(21, 6).Value = 1705
$ws.Cells.Item(22, 6).Value = 2124
$ws.Cells.Item(23, 6).Value = 565
$ws.Cells.Item(24, 6).Value = 60
$ws.Cells.Item(25, 6).Value = 1818
$ws.Cells.Item(26, 6).Value = 280
$ws.Cells.Item(27, 6).Value = 2649
$ws.Cells.Item(28, 6).Value = 474
$ws.Cells.Item(29, 6).Value = 73
$ws.Cells.Item(30, 6).Value = 665
$ws.Cells.Item(34, 6).Value = 919
$ws.Cells.Item(35, 6).Value = 1637
$ws.Cells.Item(36, 6).Value = 298
$ws.Cells.Item(38, 6).Value = 520
$ws.Cells.Item(39, 6).Value = 131
$ws.Cells.Item(40, 6).Value = 107
$ws.Cells.Item(41, 6).Value = 145

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 132
$ws.Cells.Item(9, 6).Value = 4
$ws.Cells.Item(10, 6).Value = 14
$ws.Cells.Item(11, 6).Value = 67

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 787
$ws.Cells.Item(4, 6).Value = 519
$ws.Cells.Item(5, 6).Value = 265
$ws.Cells.Item(6, 6).Value = 478
$ws.Cells.Item(7, 6).Value = 1117
$ws.Cells.Item(8, 6).Value = 316
$ws.Cells.Item(9, 6).Value = 22
$ws.Cells.Item(12, 6).Value = 1108
$ws.Cells.Item(14, 6).Value = 755
$ws.Cells.Item(15, 6).Value = 796
$ws.Cells.Item(16, 6).Value = 177
$ws.Cells.Item(17, 6).Value = 132
$ws.Cells.Item(18, 6).Value = 132
$ws.Cells.Item(22, 6).Value = 59
$ws.Cells.Item(23, 6).Value = 661
$ws.Cells.Item(24, 6).Value = 159
$ws.Cells.Item(25, 6).Value = 1705
$ws.Cells.Item(26, 6).Value = 2124
$ws.Cells.Item(27, 6).Value = 565
$ws.Cells.Item(28, 6).Value = 60
$ws.Cells.Item(31, 6).Value = 2649
$ws.Cells.Item(32, 6).Value = 474
$ws.Cells.Item(34, 6).Value = 4
$ws.Cells.Item(35, 6).Value = 14
$ws.Cells.Item(36, 6).Value = 73
$ws.Cells.Item(37, 6).Value = 67
$ws.Cells.Item(38, 6).Value = 665
$ws.Cells.Item(42, 6).Value = 919
$ws.Cells.Item(43, 6).Value = 1637
$ws.Cells.Item(45, 6).Value = 298
$ws.Cells.Item(46, 6).Value = 520
$ws.Cells.Item(47, 6).Value = 131
$ws.Cells.Item(48, 6).Value = 107
$ws.Cells.Item(49, 6).Value = 145
